$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New commit hash for the existing row 14 task ("Front userProfile" work continues)
# Copy the formatting used by the other "commit" column cells (small Consolas font)
# onto the newly filled F14 cell.
$ws.Range("F11").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F14").Value = "fd764e02ce212fa77881bf2fac85a1d6aeef3c83"

# New row 15: second work session on the same day (43984 = 2020-06-02)
$ws.Range("A15").Value = 43984
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = 0.95347222222222217
$ws.Range("E15").Value = "Zoom Integration + Front userProfile"

# Move the active selection to F16, matching where the user clicked next
$ws.Range("F16").Select()
